$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,12
$arr[0,0] = 0.6621999664170914
$arr[0,1] = 0.0555262346195704
$arr[0,2] = 0.0821304690908633
$arr[0,3] = 0
$arr[0,4] = 2.630051839096126
$arr[0,5] = 1.918020106059501
$arr[0,6] = 1.661074597085232
$arr[0,7] = 1.722858182555314
$arr[0,8] = 0.2170069821020046
$arr[0,9] = 0.4209538336900209
$arr[0,10] = 0.347835227355958
$arr[0,11] = 0.2295446117592661
$arr[1,0] = 0.6340892226893118
$arr[1,1] = 0.05237827286197927
$arr[1,2] = 0.08143863797669582
$arr[1,3] = 0
$arr[1,4] = 2.63504037662328
$arr[1,5] = 1.920970459983408
$arr[1,6] = 1.667346029514547
$arr[1,7] = 1.729665702874065
$arr[1,8] = 0.2180063437241273
$arr[1,9] = 0.3918683684832729
$arr[1,10] = 0.3462438254031426
$arr[1,11] = 0.2244553347820073
$arr[2,0] = 0.6171293087827507
$arr[2,1] = 0.05042173746649325
$arr[2,2] = 0.08103920037361689
$arr[2,3] = 0
$arr[2,4] = 2.639128511197484
$arr[2,5] = 1.923611485147148
$arr[2,6] = 1.671755179638438
$arr[2,7] = 1.734466698649129
$arr[2,8] = 0.2186626385924466
$arr[2,9] = 0.3741839674349592
$arr[2,10] = 0.345395422015649
$arr[2,11] = 0.2214325508735335
$arr[3,0] = 0.6102939948056871
$arr[3,1] = 0.04961846869080944
$arr[3,2] = 0.08088283942574037
$arr[3,3] = 0
$arr[3,4] = 2.64105248117815
$arr[3,5] = 1.924896380033658
$arr[3,6] = 1.673692522534523
$arr[3,7] = 1.736579448249032
$arr[3,8] = 0.2189408365562091
$arr[3,9] = 0.3670215991003829
$arr[3,10] = 0.345082157816563
$arr[3,11] = 0.2202265238073942
$arr[4,0] = 0.6091636004759664
$arr[4,1] = 0.04948472593314079
$arr[4,2] = 0.08085726429104767
$arr[4,3] = 0
$arr[4,4] = 2.641387545868589
$arr[4,5] = 1.925122340285583
$arr[4,6] = 1.674022711794564
$arr[4,7] = 1.736939712425425
$arr[4,8] = 0.218987681114049
$arr[4,9] = 0.3658349718185718
$arr[4,10] = 0.3450321043667159
$arr[4,11] = 0.2200278246358209
$arr[5,0] = 0.6170368170054701
$arr[5,1] = 0.0504109284683949
$arr[5,2] = 0.08103706561664836
$arr[5,3] = 0
$arr[5,4] = 2.639153413495286
$arr[5,5] = 1.923627968757742
$arr[5,6] = 1.671780737932892
$arr[5,7] = 1.734494558912004
$arr[5,8] = 0.2186663469047296
$arr[5,9] = 0.3740871939051971
$arr[5,10] = 0.3453910656357095
$arr[5,11] = 0.2214161814210982
$arr[6,0] = 0.652445366365356
$arr[6,1] = 0.0544457231430755
$arr[6,2] = 0.08188668728032411
$arr[6,3] = 0
$arr[6,4] = 2.631559255818189
$arr[6,5] = 1.918865266040243
$arr[6,6] = 1.663121164910351
$arr[6,7] = 1.725076561219176
$arr[6,8] = 0.2173427149450067
$arr[6,9] = 0.4108892213088779
$arr[6,10] = 0.3472598530923605
$arr[6,11] = 0.2277687177128271
$arr[7,0] = 0.7242455262577892
$arr[7,1] = 0.06217095319045995
$arr[7,2] = 0.08375237481870812
$arr[7,3] = 0
$arr[7,4] = 2.624792621829314
$arr[7,5] = 1.916105793324377
$arr[7,6] = 1.650565035623941
$arr[7,7] = 1.71153222742727
$arr[7,8] = 0.2150848663182021
$arr[7,9] = 0.4844279644937615
$arr[7,10] = 0.3519420855267725
$arr[7,11] = 0.2410313727644393
$arr[8,0] = 0.7784206034687315
$arr[8,1] = 0.06773433901160786
$arr[8,2] = 0.08524289510792471
$arr[8,3] = 0
$arr[8,4] = 2.624765843096256
$arr[8,5] = 1.918090310880288
$arr[8,6] = 1.644030745806702
$arr[8,7] = 1.704578536124608
$arr[8,8] = 0.2136307379488382
$arr[8,9] = 0.5392818892232754
$arr[8,10] = 0.3559980529421409
$arr[8,11] = 0.2512617993573087
$arr[9,0] = 0.8033720087977372
$arr[9,1] = 0.07024125577858342
$arr[9,2] = 0.08594660500303775
$arr[9,3] = 0
$arr[9,4] = 2.625825682248319
$arr[9,5] = 1.919864545214054
$arr[9,6] = 1.641640972551556
$arr[9,7] = 1.702065061927925
$arr[9,8] = 0.2130134123209526
$arr[9,9] = 0.5644137278142978
$arr[9,10] = 0.3579760826053615
$arr[9,11] = 0.2560206027545249
$arr[10,0] = 0.8128641745951199
$arr[10,1] = 0.07118714341790167
$arr[10,2] = 0.08621673656116968
$arr[10,3] = 0
$arr[10,4] = 2.626380987108334
$arr[10,5] = 1.920661697587576
$arr[10,6] = 1.640819694206485
$arr[10,7] = 1.701206628215203
$arr[10,8] = 0.2127859784340629
$arr[10,9] = 0.5739558907647506
$arr[10,10] = 0.3587441415987627
$arr[10,11] = 0.2578376254156041
$arr[11,0] = 0.8108179355458276
$arr[11,1] = 0.07098358188225973
$arr[11,2] = 0.0861583970911397
$arr[11,3] = 0
$arr[11,4] = 2.626254547553941
$arr[11,5] = 1.920484444768576
$arr[11,6] = 1.640992851473371
$arr[11,7] = 1.701387356060827
$arr[11,8] = 0.2128346789552751
$arr[11,9] = 0.5718996970325918
$arr[11,10] = 0.3585778818637095
$arr[11,11] = 0.2574456333446804
$arr[12,0] = 0.8041520639521309
$arr[12,1] = 0.07031914316280563
$arr[12,2] = 0.08596875591997843
$arr[12,3] = 0
$arr[12,4] = 2.625868282868424
$arr[12,5] = 1.91992761661345
$arr[12,6] = 1.641571729138008
$arr[12,7] = 1.701992567337669
$arr[12,8] = 0.2129945743526171
$arr[12,9] = 0.5651982628533858
$arr[12,10] = 0.3580388905983796
$arr[12,11] = 0.2561697910358163
$arr[13,0] = 0.8000746890451182
$arr[13,1] = 0.06991170918296064
$arr[13,2] = 0.08585306966865147
$arr[13,3] = 0
$arr[13,4] = 2.625651730025837
$arr[13,5] = 1.919602858708771
$arr[13,6] = 1.641937202204772
$arr[13,7] = 1.702375433295181
$arr[13,8] = 0.213093339266301
$arr[13,9] = 0.5610967217523637
$arr[13,10] = 0.3577112172815191
$arr[13,11] = 0.2553902463314941
$arr[14,0] = 0.7767961013422564
$arr[14,1] = 0.06757002602084583
$arr[14,2] = 0.08519741914143708
$arr[14,3] = 0
$arr[14,4] = 2.624718135790616
$arr[14,5] = 1.917991894287042
$arr[14,6] = 1.644198638768344
$arr[14,7] = 1.704755866661984
$arr[14,8] = 0.2136719688891446
$arr[14,9] = 0.5376430233288261
$arr[14,10] = 0.3558714511331544
$arr[14,11] = 0.250952901416241
$arr[15,0] = 0.7625936729053535
$arr[15,1] = 0.06612736532670738
$arr[15,2] = 0.08480174454405898
$arr[15,3] = 0
$arr[15,4] = 2.624419841970052
$arr[15,5] = 1.917226803882926
$arr[15,6] = 1.645735119722588
$arr[15,7] = 1.706382559777069
$arr[15,8] = 0.2140382389936857
$arr[15,9] = 0.5233003832999543
$arr[15,10] = 0.3547767974477125
$arr[15,11] = 0.2482575225308281
$arr[16,0] = 0.7544537287092226
$arr[16,1] = 0.06529533538611076
$arr[16,2] = 0.0845765809652903
$arr[16,3] = 0
$arr[16,4] = 2.624349186166228
$arr[16,5] = 1.916868772393727
$arr[16,6] = 1.64667371701816
$arr[16,7] = 1.707379356861985
$arr[16,8] = 0.2142530658595501
$arr[16,9] = 0.5150677017448686
$arr[16,10] = 0.3541597011730744
$arr[16,11] = 0.2467170957458009
$arr[17,0] = 0.7517026658850341
$arr[17,1] = 0.06501323769398937
$arr[17,2] = 0.08450076091927627
$arr[17,3] = 0
$arr[17,4] = 2.624342601316627
$arr[17,5] = 1.916761638363795
$arr[17,6] = 1.647000934723181
$arr[17,7] = 1.707727362892449
$arr[17,8] = 0.2143265172007265
$arr[17,9] = 0.5122831589221448
$arr[17,10] = 0.3539529161550092
$arr[17,11] = 0.2461972349879389
$arr[18,0] = 0.7641025566351516
$arr[18,1] = 0.06628117178989612
$arr[18,2] = 0.08484361480968516
$arr[18,3] = 0
$arr[18,4] = 2.624441152904538
$arr[18,5] = 1.917299759639619
$arr[18,6] = 1.645565882347128
$arr[18,7] = 1.706203065586166
$arr[18,8] = 0.2139988187023718
$arr[18,9] = 0.524825443248119
$arr[18,10] = 0.3548920301684859
$arr[18,11] = 0.2485434281566938
$arr[19,0] = 0.8061088125927256
$arr[19,1] = 0.07051439786775404
$arr[19,2] = 0.08602435932435526
$arr[19,3] = 0
$arr[19,4] = 2.625977561055919
$arr[19,5] = 1.920087770471142
$arr[19,6] = 1.641399428686825
$arr[19,7] = 1.701812268829492
$arr[19,8] = 0.2129474374318185
$arr[19,9] = 0.5671659534861249
$arr[19,10] = 0.3581966899067055
$arr[19,11] = 0.2565441316149517
$arr[20,0] = 0.833816282360516
$arr[20,1] = 0.07326108277943888
$arr[20,2] = 0.08681731238102941
$arr[20,3] = 0
$arr[20,4] = 2.627879041623373
$arr[20,5] = 1.922640160160213
$arr[20,6] = 1.639164090339065
$arr[20,7] = 1.699486793965804
$arr[20,8] = 0.2122972102337854
$arr[20,9] = 0.5949851028891828
$arr[20,10] = 0.36046730710585
$arr[20,11] = 0.2618602372969576
$arr[21,0] = 0.8190052189299877
$arr[21,1] = 0.07179694997942931
$arr[21,2] = 0.08639216507533121
$arr[21,3] = 0
$arr[21,4] = 2.626782140345426
$arr[21,5] = 1.921211091742961
$arr[21,6] = 1.640312547788497
$arr[21,7] = 1.700678177334169
$arr[21,8] = 0.2126408768970123
$arr[21,9] = 0.5801241605290102
$arr[21,10] = 0.3592453265302709
$arr[21,11] = 0.2590149944912952
$arr[22,0] = 0.7634203117014522
$arr[22,1] = 0.06621164406658409
$arr[22,2] = 0.08482467804953586
$arr[22,3] = 0
$arr[22,4] = 2.624431204109186
$arr[22,5] = 1.917266521442173
$arr[22,6] = 1.645642222436948
$arr[22,7] = 1.706284023051573
$arr[22,8] = 0.2140166273663677
$arr[22,9] = 0.5241359228336364
$arr[22,10] = 0.354839895339083
$arr[22,11] = 0.2484141416121091
$arr[23,0] = 0.7045704283137013
$arr[23,1] = 0.06010092149210777
$arr[23,2] = 0.0832264846306856
$arr[23,3] = 0
$arr[23,4] = 2.625754203362249
$arr[23,5] = 1.916147804645021
$arr[23,6] = 1.653488779369155
$arr[23,7] = 1.714669630212271
$arr[23,8] = 0.2156596372810426
$arr[23,9] = 0.4643881736899687
$arr[23,10] = 0.3505669377306759
$arr[23,11] = 0.2373577333117822

$ws.Range("B2:M25").Value = $arr

Write-Output "done"